# temp solve of RWheel
# Set the Fitness column (C) for rows 2-12 to a constant value of 4319.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C12").Value = 4319
